$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Insert 4 blank rows before row 16 to make room for the new explanatory
# text (this also shifts the existing rows 16-23 down to 20-27, and Excel
# auto-updates the cross-sheet formulas that reference them).
$ws.Rows("15:18").Insert()

$ws.Range("A15").Value = "Because enhanced rock weathering involves the crushing of rock and application"
$ws.Range("A16").Value = "on fields, which are available technologies today, we assume the full potential"
$ws.Range("A17").Value = "could be realized as soon as 2030. However, it is worth noting that would require"
$ws.Range("A18").Value = "a very large scale-up of capacity and transportation logistics."

# Widen column B to fit the new text (53 characters, stored internally with
# Excel's own padding offset, hence the 52 + 1/6 input).
$ws.Columns("B").ColumnWidth = 52 + 1/6

[void]$ws.Range("B16").Select()
